$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.Value = "'" + $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "317.13"
Set-TextValue $ws.Range("E2") "-3.58%"
Set-TextValue $ws.Range("D3") "41.71"
Set-TextValue $ws.Range("E3") "-4.94%"
Set-TextValue $ws.Range("E4") "-3.14%"
Set-TextValue $ws.Range("D5") "0.08096"
Set-TextValue $ws.Range("E5") "-3.51%"
Set-TextValue $ws.Range("D6") "4.361"
Set-TextValue $ws.Range("E6") "-1.69%"
Set-TextValue $ws.Range("D7") "1.749"
Set-TextValue $ws.Range("E7") "-10.04%"
Set-TextValue $ws.Range("D8") "0.9295"
Set-TextValue $ws.Range("E8") "-4.75%"
Set-TextValue $ws.Range("D9") "0.1115"
Set-TextValue $ws.Range("E9") "-0.98%"
Set-TextValue $ws.Range("D10") "0.1853"
Set-TextValue $ws.Range("E10") "-2.60%"
Set-TextValue $ws.Range("D11") "0.09269"
Set-TextValue $ws.Range("E11") "-4.25%"
Set-TextValue $ws.Range("D12") "0.04584"
Set-TextValue $ws.Range("E12") "-0.66%"
Set-TextValue $ws.Range("D13") "7.406"
Set-TextValue $ws.Range("E13") "-18.00%"
Set-TextValue $ws.Range("D14") "0.1055"
Set-TextValue $ws.Range("E14") "-0.75%"
Set-TextValue $ws.Range("D15") "0.001294"
Set-TextValue $ws.Range("E15") "-0.03%"
Set-TextValue $ws.Range("D16") "0.005844"
Set-TextValue $ws.Range("E16") "-4.39%"
Set-TextValue $ws.Range("E17") "-1.77%"
Set-TextValue $ws.Range("D18") "2.599"
Set-TextValue $ws.Range("E18") "3.64%"
Set-TextValue $ws.Range("D19") "0.3385"
Set-TextValue $ws.Range("E19") "1.68%"
Set-TextValue $ws.Range("D20") "0.1383"
Set-TextValue $ws.Range("E20") "0.88%"
Set-TextValue $ws.Range("D21") "0.2552"
Set-TextValue $ws.Range("E21") "0.11%"
Set-TextValue $ws.Range("D22") "0.04187"
Set-TextValue $ws.Range("E22") "0.53%"
Set-TextValue $ws.Range("D23") "0.001245"
Set-TextValue $ws.Range("E23") "-3.93%"
Set-TextValue $ws.Range("D24") "0.004245"
Set-TextValue $ws.Range("E24") "-3.60%"
Set-TextValue $ws.Range("D25") "0.0001225"
Set-TextValue $ws.Range("E25") "-5.91%"
Set-TextValue $ws.Range("E26") "0.03%"
Set-TextValue $ws.Range("D38") "0.02584"
Set-TextValue $ws.Range("E38") "-3.03%"
Set-TextValue $ws.Range("D39") "0.05474"
Set-TextValue $ws.Range("E39") "-2.99%"
Set-TextValue $ws.Range("D40") "0.008051"
Set-TextValue $ws.Range("E40") "2.39%"
Set-TextValue $ws.Range("D41") "0.1390"
Set-TextValue $ws.Range("E41") "-1.87%"
Set-TextValue $ws.Range("D42") "0.006553"
Set-TextValue $ws.Range("E42") "-10.99%"
Set-TextValue $ws.Range("D43") "0.002088"
Set-TextValue $ws.Range("E43") "-1.20%"
Set-TextValue $ws.Range("D44") "0.008237"
Set-TextValue $ws.Range("E44") "4.16%"
Set-TextValue $ws.Range("D45") "0.3450"
Set-TextValue $ws.Range("E45") "-1.97%"
Set-TextValue $ws.Range("D46") "0.00006753"
Set-TextValue $ws.Range("E46") "-2.20%"
Set-TextValue $ws.Range("E47") "0.20%"
Set-TextValue $ws.Range("D48") "0.003396"
Set-TextValue $ws.Range("E48") "-3.19%"
Set-TextValue $ws.Range("E49") "16.37%"
Set-TextValue $ws.Range("E50") "0.20%"
Set-TextValue $ws.Range("E51") "0.20%"
